$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 123, pushing the existing rows 123-140 down to 124-141.
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new weekly record.
$ws.Range("A123").Value = 7
$ws.Range("B123").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C123").Value = "Ñuble"
$ws.Range("D123").Value = 44644
$ws.Range("E123").Value = 16
$ws.Range("F123").Value = 100112045
$ws.Range("G123").Value = "Zapallo"
$ws.Range("H123").Value = "Camote"
$ws.Range("I123").Value = "1a (cosecha)"
$ws.Range("J123").Value = 200
$ws.Range("K123").Value = 300
$ws.Range("L123").Value = 350
$ws.Range("M123").Value = 325
$ws.Range("N123").Value = "$/kilo (volumen en unidades)"
$ws.Range("O123").Value = "Región de O'Higgins"
$ws.Range("P123").Value = 325
$ws.Range("Q123").Value = 1
$ws.Range("R123").Value = "Hortaliza"
